$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Give each D-column reservation a unique "owner email" value / hyperlink
# instead of several rows sharing the same address.
$ws.Range("D3").Value = "zhekaprosto7@gmail.com"
$ws.Range("D6").Value = "d.chubenko996@gmail.com"
$ws.Range("D7").Value = "khreptunchik@gmail.com"

# Rebuild the mailto hyperlinks so each one lines back up with its cell
# (order/target assignment changed as part of the same cleanup).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:yevhen@test.ca")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:max@test.ca")
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:kate@test.ca")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:brandon@test.ca")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:diana@test.ca")
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:yevhen@test.ca")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:khreptunchik@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:d.chubenko996@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:khreptunchik@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:d.chubenko996@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:khreptunchik@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:khreptunchik@gmail.com")

# Adding hyperlinks re-applies the "Hyperlink" cell style; reassert it so
# every touched cell keeps using the workbook's existing hyperlink style
# record instead of a freshly minted duplicate.
$ws.Range("G2:G7").Style = "Гиперссылка"
$ws.Range("D2:D7").Style = "Гиперссылка"

# Leave the active selection on D3, matching the saved view state.
$ws.Range("D3").Select()
